# fix pre and post time loading for FcwEventSegmenter
# Insert two new parameter rows (PRE_TIME_FCW / POST_TIME_FCW) for the
# new FcwEventDetector class, right after the existing AEB rows (row 5)
# in the "params" sheet, pushing the FcwKpiExtractor / BaseEventSegmenter
# rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Insert two blank rows at row 13 (before the WINDOW_S / FcwKpiExtractor row)
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# New row 13: PRE_TIME_FCW
$ws.Range("A13").Value = "PRE_TIME_FCW"
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = "float"
$ws.Range("D13").Value = "s"
$ws.Range("E13").Value = "time before event (duration)"
$ws.Range("F13").Value = "FcwEventDetector"

# New row 14: POST_TIME_FCW
$ws.Range("A14").Value = "POST_TIME_FCW"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "float"
$ws.Range("D14").Value = "s"
$ws.Range("E14").Value = "time after event (duration)"
$ws.Range("F14").Value = "FcwEventDetector"

# Update the selected cell to match the author's saved selection state
$ws.Range("E13").Select() | Out-Null
